$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 147; existing rows 147-224 shift down to 148-225.
$ws.Rows.Item(147).Insert()

# Populate the newly inserted row 147 with its data.
$ws.Cells.Item(147, 1).Value = 3
$ws.Cells.Item(147, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(147, 3).Value = "Coquimbo"
$ws.Cells.Item(147, 4).Value = 44518
$ws.Cells.Item(147, 5).Value = 5
$ws.Cells.Item(147, 6).Value = 100112012
$ws.Cells.Item(147, 7).Value = "Espinaca"
$ws.Cells.Item(147, 8).Value = "Sin especificar"
$ws.Cells.Item(147, 9).Value = "Primera"
$ws.Cells.Item(147, 10).Value = 160
$ws.Cells.Item(147, 11).Value = 2000
$ws.Cells.Item(147, 12).Value = 2000
$ws.Cells.Item(147, 13).Value = 2000
$ws.Cells.Item(147, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(147, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(147, 16).Value = 667
$ws.Cells.Item(147, 17).Value = 3
$ws.Cells.Item(147, 18).Value = "Hortaliza"
